$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-17 (NATMI recomputation per commit message "Natmi following Dr Hou advice")
$data = @{
    2  = @{ E=3; G=16.506869; H=49.520607; I=0.2165594803671733; J=0.2165594803671733; K=3;
            M=2.993142333333334; N=8.979427000000001; O=0.03484385887642424; P=0.03484385887642424;
            Q=49.40740839468766; R=444.666675552189; S=0.007545767972265551; T=0.007545767972265554 }
    3  = @{ E=3; G=16.506869; H=49.520607; I=0.2165594803671733; J=0.2165594803671733; K=3;
            M=31.995262; N=95.985786; O=0.3724642097459734; P=0.3724642097459735;
            Q=528.1415984546779; R=4753.274386092102; S=0.08066065571795784; T=0.08066065571795787 }
    4  = @{ E=3; G=16.506869; H=49.520607; I=0.2165594803671733; J=0.2165594803671733; K=3;
            M=34.28929533333334; N=102.867886; O=0.3991695798295478; P=0.3991695798295478;
            Q=566.0089061696447; R=5094.080155526802; S=0.08644395678626976; T=0.08644395678626977 }
    5  = @{ E=3; G=16.506869; H=49.520607; I=0.2165594803671733; J=0.2165594803671733; K=3;
            M=16.62387466666667; N=49.871624; O=0.1935223515480544; P=0.1935223515480545;
            Q=274.4081213950853; R=2469.673092555768; S=0.0419090998906801; T=0.04190909989068012 }
    6  = @{ E=3; G=24.781512; H=74.34453600000001; I=0.3251174623990092; J=0.3251174623990092; K=3;
            M=2.993142333333334; N=8.979427000000001; O=0.03484385887642424; P=0.03484385887642424;
            Q=74.17459265120802; R=667.5713338608722; S=0.01132834697809224; T=0.01132834697809224 }
    7  = @{ E=3; G=24.781512; H=74.34453600000001; I=0.3251174623990092; J=0.3251174623990092; K=3;
            M=31.995262; N=95.985786; O=0.3724642097459734; P=0.3724642097459735;
            Q=792.8909691961441; R=7136.018722765297; S=0.1210946187070632; T=0.1210946187070632 }
    8  = @{ E=3; G=24.781512; H=74.34453600000001; I=0.3251174623990092; J=0.3251174623990092; K=3;
            M=34.28929533333334; N=102.867886; O=0.3991695798295478; P=0.3991695798295478;
            Q=849.7405837745441; R=7647.665253970897; S=0.1297770008610613; T=0.1297770008610613 }
    9  = @{ E=3; G=24.781512; H=74.34453600000001; I=0.3251174623990092; J=0.3251174623990092; K=3;
            M=16.62387466666667; N=49.871624; O=0.1935223515480544; P=0.1935223515480545;
            Q=411.9647495384961; R=3707.682745846464; S=0.06291749585279242; T=0.06291749585279244 }
    10 = @{ E=3; G=24.67943933333333; H=74.038318; I=0.3237783348120013; J=0.3237783348120013; K=3;
            M=2.993142333333334; N=8.979427000000001; O=0.03484385887642424; P=0.03484385887642424;
            Q=73.86907463153179; R=664.8216716837861; S=0.01128168660543301; T=0.01128168660543301 }
    11 = @{ E=3; G=24.67943933333333; H=74.038318; I=0.3237783348120013; J=0.3237783348120013; K=3;
            M=31.995262; N=95.985786; O=0.3724642097459734; P=0.3724642097459735;
            Q=789.6251274831054; R=7106.626147347949; S=0.1205958416086193; T=0.1205958416086193 }
    12 = @{ E=3; G=24.67943933333333; H=74.038318; I=0.3237783348120013; J=0.3237783348120013; K=3;
            M=34.28929533333334; N=102.867886; O=0.3991695798295478; P=0.3991695798295478;
            Q=846.2405839617498; R=7616.165255655748; S=0.1292424618648172; T=0.1292424618648172 }
    13 = @{ E=3; G=24.67943933333333; H=74.038318; I=0.3237783348120013; J=0.3237783348120013; K=3;
            M=16.62387466666667; N=49.871624; O=0.1935223515480544; P=0.1935223515480545;
            Q=410.2679063209369; R=3692.411156888432; S=0.0626583447331318; T=0.0626583447331318 }
    14 = @{ E=3; G=10.255437; H=30.766311; I=0.1345447224218162; J=0.1345447224218162; K=3;
            M=2.993142333333334; N=8.979427000000001; O=0.03484385887642424; P=0.03484385887642424;
            Q=30.695982631533; R=276.2638436837971; S=0.004688057320633435; T=0.004688057320633437 }
    15 = @{ E=3; G=10.255437; H=30.766311; I=0.1345447224218162; J=0.1345447224218162; K=3;
            M=31.995262; N=95.985786; O=0.3724642097459734; P=0.3724642097459735;
            Q=328.125393739494; R=2953.128543655446; S=0.05011309371233313; T=0.05011309371233314 }
    16 = @{ E=3; G=10.255437; H=30.766311; I=0.1345447224218162; J=0.1345447224218162; K=3;
            M=34.28929533333334; N=102.867886; O=0.3991695798295478; P=0.3991695798295478;
            Q=351.651708065394; R=3164.865372588546; S=0.05370616031739951; T=0.05370616031739952 }
    17 = @{ E=3; G=10.255437; H=30.766311; I=0.1345447224218162; J=0.1345447224218162; K=3;
            M=16.62387466666667; N=49.871624; O=0.1935223515480544; P=0.1935223515480545;
            Q=170.485099339896; R=1534.365894059064; S=0.02603741107145012; T=0.02603741107145013 }
}

foreach ($rowNum in $data.Keys) {
    $rowVals = $data[$rowNum]
    foreach ($col in $rowVals.Keys) {
        $ws.Range("$col$rowNum").Value = $rowVals[$col]
    }
}
